$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the "(according to the population census data)" subtitle row)
$ws.Rows.Item(2).Delete()

# Delete columns B:C (drop the 1989 and 2002 data columns, keep only 2014)
$ws.Range("B1:C1").EntireColumn.Delete()

# Re-apply row heights: rows 1-6 get a custom height of 20.1
$ws.Range("A1:B6").RowHeight = 20.1
